# Applies updated match-stats values (Jogos_totais, Jogos_titular, minutos_partida,
# minutos_totais, gols, assist) for the affected player rows in the
# brasileirão_2024 sheet, per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 14
$ws.Range("G2").Value = 1260

# Row 6
$ws.Range("D6").Value = 24
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 68
$ws.Range("G6").Value = 1624
$ws.Range("H6").Value = 1

# Row 8
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 88
$ws.Range("G8").Value = 963

# Row 10
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 58
$ws.Range("G10").Value = 871

# Row 11
$ws.Range("D11").Value = 24
$ws.Range("E11").Value = 24
$ws.Range("G11").Value = 2145
$ws.Range("H11").Value = 4

# Row 15
$ws.Range("D15").Value = 28
$ws.Range("E15").Value = 28
$ws.Range("G15").Value = 2373
$ws.Range("H15").Value = 6
$ws.Range("I15").Value = 5

# Row 16
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 244

# Row 17
$ws.Range("D17").Value = 25
$ws.Range("F17").Value = 49
$ws.Range("G17").Value = 1220

# Row 18
$ws.Range("D18").Value = 25
$ws.Range("F18").Value = 71
$ws.Range("G18").Value = 1780

# Row 19
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 61
$ws.Range("G19").Value = 1283

# Row 20
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 6
$ws.Range("F20").Value = 75
$ws.Range("G20").Value = 448

# Row 24
$ws.Range("D24").Value = 9
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 492

# Row 27
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 54
$ws.Range("G27").Value = 214
$ws.Range("H27").Value = 1

# Row 28
$ws.Range("D28").Value = 23
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 77
$ws.Range("G28").Value = 1761
$ws.Range("H28").Value = 8

# Row 29
$ws.Range("D29").Value = 21
$ws.Range("G29").Value = 1313

# Row 30
$ws.Range("D30").Value = 7
$ws.Range("F30").Value = 32
$ws.Range("G30").Value = 227
$ws.Range("I30").Value = 1

# Row 131
$ws.Range("D131").Value = 27
$ws.Range("E131").Value = 27
$ws.Range("G131").Value = 2430

# Row 134
$ws.Range("D134").Value = 3
$ws.Range("F134").Value = 46
$ws.Range("G134").Value = 138

# Row 135
$ws.Range("D135").Value = 25
$ws.Range("E135").Value = 22
$ws.Range("G135").Value = 1840

# Row 137
$ws.Range("D137").Value = 22
$ws.Range("E137").Value = 21
$ws.Range("G137").Value = 1857

# Row 138
$ws.Range("D138").Value = 23
$ws.Range("E138").Value = 17
$ws.Range("F138").Value = 72
$ws.Range("G138").Value = 1648

# Row 140
$ws.Range("D140").Value = 24
$ws.Range("E140").Value = 19
$ws.Range("F140").Value = 75
$ws.Range("G140").Value = 1794

# Row 141
$ws.Range("D141").Value = 20
$ws.Range("E141").Value = 16
$ws.Range("G141").Value = 1432

# Row 143
$ws.Range("D143").Value = 14
$ws.Range("F143").Value = 66
$ws.Range("G143").Value = 918

# Row 144
$ws.Range("D144").Value = 21
$ws.Range("E144").Value = 16
$ws.Range("G144").Value = 1409

# Row 148
$ws.Range("D148").Value = 3
$ws.Range("F148").Value = 44
$ws.Range("G148").Value = 133

# Row 150
$ws.Range("D150").Value = 5
$ws.Range("E150").Value = 3
$ws.Range("F150").Value = 63
$ws.Range("G150").Value = 317

# Row 151
$ws.Range("D151").Value = 4
$ws.Range("F151").Value = 49
$ws.Range("G151").Value = 194

# Row 153
$ws.Range("D153").Value = 17
$ws.Range("F153").Value = 41
$ws.Range("G153").Value = 698

# Row 154
$ws.Range("D154").Value = 22
$ws.Range("E154").Value = 15
$ws.Range("F154").Value = 60
$ws.Range("G154").Value = 1322

# Row 157
$ws.Range("D157").Value = 15
$ws.Range("E157").Value = 4
$ws.Range("F157").Value = 37
$ws.Range("G157").Value = 555

# Row 158
$ws.Range("D158").Value = 9
$ws.Range("E158").Value = 5
$ws.Range("F158").Value = 52
$ws.Range("G158").Value = 464

# Row 484
$ws.Range("D484").Value = 9
$ws.Range("E484").Value = 7
$ws.Range("F484").Value = 74
$ws.Range("G484").Value = 663

# Row 486
$ws.Range("D486").Value = 24
$ws.Range("E486").Value = 22
$ws.Range("G486").Value = 1950
$ws.Range("I486").Value = 4

# Row 488
$ws.Range("D488").Value = 18
$ws.Range("E488").Value = 18
$ws.Range("F488").Value = 88
$ws.Range("G488").Value = 1575

# Row 489
$ws.Range("D489").Value = 16
$ws.Range("E489").Value = 14
$ws.Range("G489").Value = 1286

# Row 495
$ws.Range("D495").Value = 1
$ws.Range("E495").Value = 1
$ws.Range("F495").Value = 79
$ws.Range("G495").Value = 79

# Row 497
$ws.Range("D497").Value = 17
$ws.Range("E497").Value = 13
$ws.Range("F497").Value = 63
$ws.Range("G497").Value = 1071

# Row 498
$ws.Range("D498").Value = 25
$ws.Range("F498").Value = 49
$ws.Range("G498").Value = 1225

# Row 500
$ws.Range("D500").Value = 17
$ws.Range("E500").Value = 11
$ws.Range("F500").Value = 58
$ws.Range("G500").Value = 986
$ws.Range("H500").Value = 3
$ws.Range("I500").Value = 2

# Row 501
$ws.Range("D501").Value = 23
$ws.Range("E501").Value = 15
$ws.Range("F501").Value = 55
$ws.Range("G501").Value = 1261

# Row 502
$ws.Range("D502").Value = 14
$ws.Range("E502").Value = 12
$ws.Range("F502").Value = 76
$ws.Range("G502").Value = 1063

# Row 503
$ws.Range("D503").Value = 26
$ws.Range("E503").Value = 23
$ws.Range("F503").Value = 76
$ws.Range("G503").Value = 1976
$ws.Range("H503").Value = 2

# Row 506
$ws.Range("D506").Value = 10
$ws.Range("F506").Value = 48
$ws.Range("G506").Value = 483

# Row 509
$ws.Range("D509").Value = 19
$ws.Range("F509").Value = 47
$ws.Range("G509").Value = 892

# Row 510
$ws.Range("D510").Value = 20
$ws.Range("E510").Value = 11
$ws.Range("F510").Value = 53
$ws.Range("G510").Value = 1059

# Row 511
$ws.Range("D511").Value = 18
$ws.Range("F511").Value = 42
$ws.Range("G511").Value = 756

# Row 512
$ws.Range("D512").Value = 27
$ws.Range("G512").Value = 1390

# Row 513
$ws.Range("D513").Value = 29
$ws.Range("E513").Value = 29
$ws.Range("G513").Value = 2610

# Row 518
$ws.Range("D518").Value = 23
$ws.Range("E518").Value = 23
$ws.Range("G518").Value = 1959

# Row 519
$ws.Range("D519").Value = 18
$ws.Range("E519").Value = 13
$ws.Range("F519").Value = 69
$ws.Range("G519").Value = 1234

# Row 520
$ws.Range("D520").Value = 19
$ws.Range("E520").Value = 16
$ws.Range("F520").Value = 76
$ws.Range("G520").Value = 1437

# Row 522
$ws.Range("D522").Value = 10
$ws.Range("E522").Value = 9
$ws.Range("F522").Value = 77
$ws.Range("G522").Value = 765

# Row 526
$ws.Range("D526").Value = 19
$ws.Range("F526").Value = 87
$ws.Range("G526").Value = 1658
$ws.Range("H526").Value = 5

# Row 527
$ws.Range("D527").Value = 17
$ws.Range("E527").Value = 15
$ws.Range("F527").Value = 71
$ws.Range("G527").Value = 1212

# Row 528
$ws.Range("D528").Value = 13
$ws.Range("F528").Value = 25
$ws.Range("G528").Value = 329

# Row 529
$ws.Range("D529").Value = 26
$ws.Range("E529").Value = 24
$ws.Range("G529").Value = 1867

# Row 531
$ws.Range("D531").Value = 25
$ws.Range("E531").Value = 25
$ws.Range("G531").Value = 1977
$ws.Range("I531").Value = 1

# Row 534
$ws.Range("D534").Value = 21
$ws.Range("E534").Value = 12
$ws.Range("F534").Value = 55
$ws.Range("G534").Value = 1159
$ws.Range("H534").Value = 3

# Row 536
$ws.Range("D536").Value = 8
$ws.Range("E536").Value = 6
$ws.Range("F536").Value = 67
$ws.Range("G536").Value = 537

# Row 538
$ws.Range("D538").Value = 18
$ws.Range("F538").Value = 71
$ws.Range("G538").Value = 1276

# Row 540
$ws.Range("D540").Value = 22
$ws.Range("E540").Value = 10
$ws.Range("F540").Value = 47
$ws.Range("G540").Value = 1032
$ws.Range("I540").Value = 2
